$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J, styled like the other headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-14
$data = @(
    @(6, 6),
    @(5, 5),
    @(9, 9),
    @(5, 5),
    @(5, 5),
    @(7, 7),
    @(7, 8),
    @(6, 6),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
